# "new changes in homepage" - add a "Message" header column (D1) to the
# Login worksheet's credential table, matching the bold header style
# already used for A1:C1 (Username / Password / Logo).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

$ws.Range("D1").Value = "Message"

# Match the existing header formatting (bold font, same as A1:C1) without
# disturbing any other style definitions.
$ws.Range("D1").Font.Bold = $true
